$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from SCD0203 to SCD0011
$ws.Name = "SCD0011"

# Update TC_ID value in B2 from "DGS-218" to "SCD0011-034"
$ws.Range("B2").Value = "SCD0011-034"

# Widen column B to fit new text (target stored width 12.42578125;
# the host's ColumnWidth setter quantizes to 1/6-character steps, so
# 11.71 is the nearest input that lands on the closest achievable width)
$ws.Columns("B").ColumnWidth = 11.71

# Update the active selection to B3
$ws.Range("B3").Select()
